$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56: normalize politeness_score (column B) to a real number (it was
# stored as text "3" before this edit).
$ws.Cells.Item(56, 2).Value = 3

# Row 57: new annotation row, appended after the (now renumbered) row 56.
$ws.Cells.Item(57, 1).Value = "Ruilin"

# politeness_score on this new row stays a text "3" (matches how the
# annotation tool originally wrote it), so force text formatting, assign
# the value, then drop back to the Normal style so no extra style index
# lingers on the cell.
$c = $ws.Cells.Item(57, 2)
$c.NumberFormat = "@"
$c.Value = "3"
$c.Style = "Normal"

$ws.Cells.Item(57, 3).Value = "无"
$ws.Cells.Item(57, 4).Value = "DIS"
$ws.Cells.Item(57, 5).Value = "WRI"
$ws.Cells.Item(57, 6).Value = "a903e5ac-dd38-46eb-9fca-86d5e31ee0d2"
$ws.Cells.Item(57, 7).Value = "H1aIuk-RW_annotated.xlsx"
$ws.Cells.Item(57, 8).Value = '2) The "Active learning" approach is simply the classing hitting set approach for computing k-center.'
